# card_export.xlsx edit: add a new "Tour" card as the first data row.
#
# Implémentation des tours : on insère une nouvelle ligne en tête du
# tableau (ligne 1) pour la nouvelle carte "Tour" (Legendaire / Royal,
# cout 10, attaque 0, vie 15, portee 0, capacites 0/1). Toutes les
# lignes existantes sont decalees d'une ligne vers le bas, ce qui
# correspond au comportement observe dans le diff (les identifiants en
# colonne A ne sont pas renumerotes, ils gardent leur valeur d'origine).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 1 - this shifts every
# existing row (and its formatting) down by one, exactly like the diff.
$ws.Rows.Item(1).Insert()

# Populate the freshly inserted row 1 with the new "Tour" card.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = "Legendaire"
$ws.Range("C1").Value = "Tour"
$ws.Range("D1").Value = "Royal"
$ws.Range("E1").Value = 10
$ws.Range("F1").Value = 0
$ws.Range("G1").Value = 15
$ws.Range("H1").Value = 0
$ws.Range("I1").Value = 0
$ws.Range("J1").Value = 1

# Row 1 uses the default (unstyled) look - give it the taller, default
# row height instead of the compact 13.2pt used by the rest of the data.
$ws.Rows.Item(1).RowHeight = 15.75

# Update the sheet view: drop the old scrolled/selected position and
# select L6 instead (also resets the scrolled top-left cell back to A1).
$ws.Activate()
$ws.Range("L6").Select()
